$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rename the bookmark "if-and-in" -> "if-in"
#    (Bookmark.Name is read-only, so delete + re-add over the same range)
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("if-and-in")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("if-in", $bmRange)

# ---------------------------------------------------------------------------
# 2. Heading3 paragraph: "if and in" -> "if in" (single VerbatimChar run)
#    Delete " and " from between "if" and "in" (merges the two pre-existing
#    VerbatimChar runs into one), then re-insert a single space so the
#    result reads "if in" instead of "ifin".
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1Search = $p1.Range
$found = $p1Search.Find.Execute("if and in", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $ifAndInStart = $p1Search.Start
    $ifAndInEnd = $p1Search.End
    $midStart = $ifAndInStart + 2
    $midEnd = $ifAndInEnd - 2
    $midRange = $d.Range($midStart, $midEnd)
    $midRange.Text = ""
    $spacePoint = $d.Range($midStart, $midStart)
    $spacePoint.InsertAfter(" ")
}

# ---------------------------------------------------------------------------
# 3. SourceCode paragraph: insert a brand-new first line followed by a line
#    break, ahead of the existing first source line.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$insertPoint = $d.Range($p3.Range.Start, $p3.Range.Start)
$newLine = '      . post `postname'' ("Variable") ("") ("Summary 1")   ("Summary 2") ("Overall")'
$insertPoint.InsertAfter($newLine)
$breakPoint = $d.Range($p3.Range.Start + $newLine.Length, $p3.Range.Start + $newLine.Length)
$breakPoint.InsertBreak(6)

# ---------------------------------------------------------------------------
# 4. Last source line ("ethnicity in 1/100  if ethnicity ==4, ...") drops
#    the " gap(2)" option - remove just the "gap(2)" token, leaving the
#    surrounding spaces (so "cat_col gap(2) n_analysis" -> "cat_col  n_analysis").
# ---------------------------------------------------------------------------
$searchText = "cat_levels(4 3 2 1 0) cat_col gap(2) n_analysis(append)"
$scan = $d.Content
$lastStart = -1
$lastEnd = -1
while ($scan.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $lastStart = $scan.Start
    $lastEnd = $scan.End
    $scan.Collapse(0)
}
if ($lastStart -ge 0) {
    $relOffset = "cat_levels(4 3 2 1 0) cat_col ".Length
    $gapLen = "gap(2)".Length
    $gapStart = $lastStart + $relOffset
    $gapEnd = $gapStart + $gapLen
    $gapRange = $d.Range($gapStart, $gapEnd)
    $gapRange.Text = ""
}
